# Add the new "2022-Q4" quarterly sheet right after "总计", and record it
# in the summary ("总计") sheet as the newest row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert a new row into the "总计" summary sheet for 2022-Q4 totals.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()

$summary.Range("B2:D2").ClearFormats()
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 0.13

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $afterSheet)
$newSheet.Name = "2022-Q4"

$refSheet = $wb.Worksheets.Item("2022-Q3")
$refSheet.Range("A1:H2").Copy($newSheet.Range("A1"))
$refSheet.Range("A2:H2").Copy($newSheet.Range("A3:A7"))
$excel.CutCopyMode = $false

# Row 2 - 002863
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'002863"
$newSheet.Range("C2").Value = "金信深圳成长灵活配置混合"
$newSheet.Range("D2").Value = "'0.73"
$newSheet.Range("E2").Value = "'91.56"
$newSheet.Range("F2").Value = "'5.31"
$newSheet.Range("G2").Value = "'0.0388"
$newSheet.Range("H2").Value = 3

# Row 3 - 007254
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'007254"
$newSheet.Range("C3").Value = "广发均衡价值混合"
$newSheet.Range("D3").Value = "'0.73"
$newSheet.Range("E3").Value = "'92.28"
$newSheet.Range("F3").Value = "'4.82"
$newSheet.Range("G3").Value = "'0.0352"
$newSheet.Range("H3").Value = 8

# Row 4 - 003142
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'003142"
$newSheet.Range("C4").Value = "鹏华弘达灵活配置混合A"
$newSheet.Range("D4").Value = "'1.18"
$newSheet.Range("E4").Value = "'39.70"
$newSheet.Range("F4").Value = "'1.84"
$newSheet.Range("G4").Value = "'0.0217"
$newSheet.Range("H4").Value = 9

# Row 5 - 001326
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'001326"
$newSheet.Range("C5").Value = "鹏华弘和灵活配置混合C"
$newSheet.Range("D5").Value = "'0.48"
$newSheet.Range("E5").Value = "'57.67"
$newSheet.Range("F5").Value = "'3.17"
$newSheet.Range("G5").Value = "'0.0152"
$newSheet.Range("H5").Value = 9

# Row 6 - 001325
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "'001325"
$newSheet.Range("C6").Value = "鹏华弘和灵活配置混合A"
$newSheet.Range("D6").Value = "'0.42"
$newSheet.Range("E6").Value = "'57.67"
$newSheet.Range("F6").Value = "'3.17"
$newSheet.Range("G6").Value = "'0.0133"
$newSheet.Range("H6").Value = 9

# Row 7 - 003143
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "'003143"
$newSheet.Range("C7").Value = "鹏华弘达灵活配置混合C"
$newSheet.Range("D7").Value = "'0.11"
$newSheet.Range("E7").Value = "'39.70"
$newSheet.Range("F7").Value = "'1.84"
$newSheet.Range("G7").Value = "'0.0020"
$newSheet.Range("H7").Value = 9
